{"js": "// Update the division problems in the practice-sheet table.\n//\n// The replacements must be applied positionally, in document order, because\n// a couple of the original expressions (e.g. \"46\u00f74=\") appear more than once\n// and each occurrence needs to become a different new value - a blind\n// search/replace across the whole document would incorrectly set every\n// occurrence to the same new text.\n\nconst replacements = [\n  [\"58\u00f75=\", \"80\u00f73=\"],\n  [\"61\u00f75=\", \"10\u00f74=\"],\n  [\"16\u00f76=\", \"50\u00f78=\"],\n  [\"37\u00f76=\", \"26\u00f73=\"],\n  [\"96\u00f74=\", \"89\u00f75=\"],\n  [\"50\u00f76=\", \"77\u00f72=\"],\n  [\"62\u00f76=\", \"19\u00f72=\"],\n  [\"51\u00f76=\", \"93\u00f78=\"],\n  [\"70\u00f76=\", \"84\u00f76=\"],\n  [\"62\u00f79=\", \"89\u00f72=\"],\n  [\"46\u00f74=\", \"96\u00f73=\"],\n  [\"97\u00f72=\", \"43\u00f73=\"],\n  [\"24\u00f73=\", \"52\u00f77=\"],\n  [\"78\u00f75=\", \"53\u00f75=\"],\n  [\"16\u00f75=\", \"11\u00f73=\"],\n  [\"46\u00f77=\", \"18\u00f77=\"],\n  [\"57\u00f75=\", \"56\u00f75=\"],\n  [\"25\u00f75=\", \"48\u00f76=\"],\n  [\"14\u00f75=\", \"39\u00f75=\"],\n  [\"46\u00f74=\", \"87\u00f79=\"],\n  [\"99\u00f72=\", \"20\u00f77=\"],\n  [\"58\u00f76=\", \"54\u00f75=\"],\n  [\"88\u00f77=\", \"77\u00f77=\"],\n  [\"59\u00f72=\", \"65\u00f76=\"],\n  [\"27\u00f79=\", \"26\u00f79=\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\nlet repIndex = 0;\nfor (let r = 0; r < table.values.length && repIndex < replacements.length; r++) {\n  const row = table.values[r];\n  for (let c = 0; c < row.length && repIndex < replacements.length; c++) {\n    const text = row[c];\n    if (!text) {\n      continue;\n    }\n\n    const [oldText, newText] = replacements[repIndex];\n    if (text === oldText) {\n      table.getCell(r, c).value = newText;\n    }\n    repIndex++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division problems in the practice-sheet table.\n#\n# The replacements must be applied positionally, in document order, because\n# a couple of the original expressions (e.g. \"46\u00f74=\") appear more than once\n# and each occurrence needs to become a different new value - a blind\n# Find/Replace-All across the whole document would incorrectly set every\n# occurrence to the same new text.\n#\n# Note: iterating with the flattened `Table.Range.Cells` collection and\n# mutating cell text while walking it leaves later items stale, so cells are\n# addressed directly through `Table.Cell(row, column)` instead.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"58\u00f75=\", \"80\u00f73=\"),\n    @(\"61\u00f75=\", \"10\u00f74=\"),\n    @(\"16\u00f76=\", \"50\u00f78=\"),\n    @(\"37\u00f76=\", \"26\u00f73=\"),\n    @(\"96\u00f74=\", \"89\u00f75=\"),\n    @(\"50\u00f76=\", \"77\u00f72=\"),\n    @(\"62\u00f76=\", \"19\u00f72=\"),\n    @(\"51\u00f76=\", \"93\u00f78=\"),\n    @(\"70\u00f76=\", \"84\u00f76=\"),\n    @(\"62\u00f79=\", \"89\u00f72=\"),\n    @(\"46\u00f74=\", \"96\u00f73=\"),\n    @(\"97\u00f72=\", \"43\u00f73=\"),\n    @(\"24\u00f73=\", \"52\u00f77=\"),\n    @(\"78\u00f75=\", \"53\u00f75=\"),\n    @(\"16\u00f75=\", \"11\u00f73=\"),\n    @(\"46\u00f77=\", \"18\u00f77=\"),\n    @(\"57\u00f75=\", \"56\u00f75=\"),\n    @(\"25\u00f75=\", \"48\u00f76=\"),\n    @(\"14\u00f75=\", \"39\u00f75=\"),\n    @(\"46\u00f74=\", \"87\u00f79=\"),\n    @(\"99\u00f72=\", \"20\u00f77=\"),\n    @(\"58\u00f76=\", \"54\u00f75=\"),\n    @(\"88\u00f77=\", \"77\u00f77=\"),\n    @(\"59\u00f72=\", \"65\u00f76=\"),\n    @(\"27\u00f79=\", \"26\u00f79=\"),\n)\n\n$t = $d.Tables(1)\n$rowCount = $t.Rows.Count\n$colCount = $t.Columns.Count\n\n$repIndex = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    if ($repIndex -ge $replacements.Count) {\n        break\n    }\n    for ($c = 1; $c -le $colCount; $c++) {\n        if ($repIndex -ge $replacements.Count) {\n            break\n        }\n\n        $cell = $t.Cell($r, $c)\n        $range = $cell.Range\n        $text = $range.Text\n        # Cell text ends with the cell-mark (chr 13 + chr 7); strip it off.\n        $text = $text.TrimEnd([char]7).TrimEnd([char]13)\n\n        if ($text.Length -eq 0) {\n            continue\n        }\n\n        $pair = $replacements[$repIndex]\n        $old = $pair[0]\n        $new = $pair[1]\n\n        if ($text -eq $old) {\n            $range.Text = $new\n        }\n\n        $repIndex++\n    }\n}\n"}
